$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.244.24'
$ws.Range("E2").Value = '  -3.99%  '
$ws.Range("D3").Value = '2.453.96'
$ws.Range("E3").Value = '  -6.87%  '
$ws.Range("D5").Value = "'548.95"
$ws.Range("D6").Value = "'146.06"
$ws.Range("E6").Value = '  -6.86%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -6.55%  '
$ws.Range("D9").Value = '2.451.50'
$ws.Range("E9").Value = '  -6.87%  '
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = '  -9.76%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = "'5.44"
$ws.Range("E11").Value = '  -6.83%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'0.154"
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = '  -7.94%  '
$ws.Range("D14").Value = "'26.04"
$ws.Range("E14").Value = '  -9.27%  '
$ws.Range("D15").Value = '2.895.07'
$ws.Range("E15").Value = '  -6.92%  '
$ws.Range("D16").Value = "'0.0000166"
$ws.Range("E16").Value = '  -10.40%  '
$ws.Range("D17").Value = '61.155.28'
$ws.Range("E17").Value = '  -4.03%  '
$ws.Range("D18").Value = '2.457.53'
$ws.Range("E18").Value = '  -6.88%  '
$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = '  -8.96%  '
$ws.Range("D20").Value = "'7.12"
$ws.Range("E20").Value = '  -8.05%  '
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = '  -7.61%  '
$ws.Range("D22").Value = "'317.92"
$ws.Range("E22").Value = '  -7.71%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = "'1.88"
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = "'63.77"
$ws.Range("E25").Value = '  -6.65%  '
$ws.Range("D26").Value = '0.0₃0978'
$ws.Range("E26").Value = '  -13.37%  '
$ws.Range("D27").Value = '2.586.30'
$ws.Range("E27").Value = '  -6.61%  '
$ws.Range("D28").Value = "'547.54"
$ws.Range("E28").Value = '  -5.54%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  -10.46%  '
$ws.Range("D31").Value = "'8.27"
$ws.Range("E31").Value = '  -10.79%  '
$ws.Range("E32").Value = '  -7.26%  '
$ws.Range("E33").Value = '  -8.98%  '
$ws.Range("E34").Value = '  -7.65%  '
$ws.Range("E35").Value = '  -8.79%  '
$ws.Range("E36").Value = '  -11.53%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -11.48%  '
$ws.Range("E39").Value = '  -6.12%  '
$ws.Range("D40").Value = "'18.35"
$ws.Range("E40").Value = '  -7.28%  '
$ws.Range("D41").Value = "'142.04"
$ws.Range("E41").Value = '  -7.00%  '
$ws.Range("D42").Value = "'1.75"
$ws.Range("E42").Value = '  -8.35%  '
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -4.20%  '
$ws.Range("D45").Value = "'2.38"
$ws.Range("E45").Value = '  -7.15%  '
$ws.Range("D46").Value = "'145.86"
$ws.Range("E46").Value = '  -10.19%  '
$ws.Range("E47").Value = '  -8.39%  '
$ws.Range("D48").Value = "'21.40"
$ws.Range("E48").Value = '  -11.59%  '
$ws.Range("E49").Value = '  -9.19%  '
$ws.Range("D50").Value = "'0.587"
$ws.Range("E50").Value = '  -7.30%  '
$ws.Range("E51").Value = '  -6.96%  '

$resetCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D31", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D50")
foreach ($c in $resetCells) {
    $ws.Range($c).Style = "Normal"
}
